# Bugfixed QoQ Visualizations and a typo in the evaluation objects
#
# The first 16 data rows of the revision_qoq_GDP sheet (quarters from
# 1984-07-01 through 1988-04-01) contained near-zero floating point
# noise values and should not have been part of the evaluation series.
# Remove them so the series now starts with the 1988-07-01 quarter,
# shifting every following row up by 16 and shrinking the used range
# from A1:B164 down to A1:B148.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-17 (the first 16 quarters after the header row) are removed;
# EntireRow.Delete() shifts all subsequent rows up automatically.
$ws.Range("A2:B17").EntireRow.Delete()
